$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-set NumberFormat to Text for cells whose new values look numeric,
# so Excel stores them as text (matching original inlineStr string data).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "34.228.76"
$ws.Range("E2").Value = "  -0.48%  "
$ws.Range("D3").Value = "1.809.04"
$ws.Range("E3").Value = "  +1.10%  "
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "224.71"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").Value = "0.554"
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").Value = "32.12"
$ws.Range("E8").Value = "  -2.30%  "
$ws.Range("D9").Value = "0.292"
$ws.Range("E9").Value = "  +3.41%  "
$ws.Range("E10").Value = "  +8.59%  "
$ws.Range("D11").Value = "0.0929"
$ws.Range("E11").Value = "  -0.27%  "
$ws.Range("D12").Value = "2.069.38"
$ws.Range("E12").Value = "  +1.00%  "
$ws.Range("D13").Value = "1.796.51"
$ws.Range("E13").Value = "  +0.60%  "
$ws.Range("D14").Value = "10.89"
$ws.Range("E14").Value = "  -1.30%  "
$ws.Range("D15").Value = "0.636"
$ws.Range("E15").Value = "  -0.65%  "
$ws.Range("D16").Value = "34.175.51"
$ws.Range("E16").Value = "  -0.71%  "
$ws.Range("E17").Value = "  +0.97%  "
$ws.Range("D18").Value = "69.42"
$ws.Range("E18").Value = "  +0.34%  "
$ws.Range("D19").Value = "248.57"
$ws.Range("E19").Value = "  -2.53%  "
$ws.Range("D20").Value = "0.0₃0794"
$ws.Range("E20").Value = "  +6.33%  "
$ws.Range("D21").Value = "10.95"
$ws.Range("E21").Value = "  +5.20%  "
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("D23").Value = "4.23"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("E24").Value = "  +0.85%  "
$ws.Range("D25").Value = "159.60"
$ws.Range("E25").Value = "  +1.04%  "
$ws.Range("D26").Value = "16.63"
$ws.Range("E26").Value = "  +1.23%  "
$ws.Range("E27").Value = "  +1.65%  "
$ws.Range("E28").Value = "  +0.22%  "
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("D30").Value = "0.0530"
$ws.Range("E30").Value = "  +3.37%  "
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("E32").Value = "  +1.58%  "
$ws.Range("E33").Value = "  -0.16%  "
$ws.Range("D34").Value = "1.87"
$ws.Range("E34").Value = "  -1.49%  "
$ws.Range("D35").Value = "1.428.27"
$ws.Range("E35").Value = "  -1.72%  "
$ws.Range("E36").Value = "  +0.80%  "
$ws.Range("E37").Value = "  +0.92%  "
$ws.Range("E38").Value = "  -0.51%  "
$ws.Range("E39").Value = "  +7.13%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "2.79"
$ws.Range("E40").Value = "  -2.11%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "80.83"
$ws.Range("E41").Value = "  -2.85%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("E43").Value = "  +3.85%  "
$ws.Range("D44").Value = "5.95"
$ws.Range("E44").Value = "  +1.08%  "
$ws.Range("D45").Value = "0.0497"
$ws.Range("E45").Value = "  -1.96%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").Value = "1.967.06"
$ws.Range("E47").Value = "  +0.77%  "
$ws.Range("D48").Value = "106.14"
$ws.Range("E48").Value = "  +6.95%  "
$ws.Range("D49").Value = "0.997"
$ws.Range("E49").Value = "  -0.35%  "
$ws.Range("D50").Value = "11.88"
$ws.Range("E50").Value = "  -3.29%  "
$ws.Range("D51").Value = "0.0₆0124"
$ws.Range("E51").Value = "  +6.10%  "

